# Update the "Datos actualizados" timestamp and refresh several countries'
# COVID-19 statistics in the "Pais" worksheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Pais")

# Update the "last updated" banner text in cell A1.
$ws.Range("A1").Value = "Datos actualizados a 13 de Octubre de 2020 a las 09:32"

# Row 5 - India
$ws.Cells.Item(5, 2).Value = 7177783   # B5 Casos totales
$ws.Cells.Item(5, 3).Value = 4218      # C5 Nuevos casos
$ws.Cells.Item(5, 5).Value = 840571    # E5 Recuperados
$ws.Cells.Item(5, 7).Value = 23        # G5 Muertes hoy
$ws.Cells.Item(5, 8).Value = 109917    # H5 Muertes

# Row 57 - Barein
$ws.Cells.Item(57, 5).Value = 3979     # E57 Recuperados
$ws.Cells.Item(57, 7).Value = 2        # G57 Muertes hoy
$ws.Cells.Item(57, 8).Value = 282      # H57 Muertes

# Row 58 - Suiza
$ws.Cells.Item(58, 5).Value = 12839    # E58 Recuperados
$ws.Cells.Item(58, 8).Value = 2097     # H58 Muertes

# Row 60 - Uzbekistan
$ws.Cells.Item(60, 2).Value = 61534    # B60 Casos totales
$ws.Cells.Item(60, 3).Value = 215      # C60 Nuevos casos
$ws.Cells.Item(60, 5).Value = 2597     # E60 Recuperados
$ws.Cells.Item(60, 7).Value = 1        # G60 Muertes hoy
$ws.Cells.Item(60, 8).Value = 510      # H60 Muertes

# Row 63 - Armenia
$ws.Cells.Item(63, 2).Value = 57566    # B63 Casos totales
$ws.Cells.Item(63, 3).Value = 745      # C63 Nuevos casos
$ws.Cells.Item(63, 4).Value = 46318    # D63 Casos activos
$ws.Cells.Item(63, 5).Value = 10216    # E63 Recuperados
$ws.Cells.Item(63, 7).Value = 6        # G63 Muertes hoy
$ws.Cells.Item(63, 8).Value = 1032     # H63 Muertes

# Row 75 - Afganistan
$ws.Cells.Item(75, 2).Value = 39928    # B75 Casos totales
$ws.Cells.Item(75, 3).Value = 58       # C75 Nuevos casos
$ws.Cells.Item(75, 4).Value = 33308    # D75 Casos activos
$ws.Cells.Item(75, 5).Value = 5140     # E75 Recuperados
$ws.Cells.Item(75, 7).Value = 1        # G75 Muertes hoy
$ws.Cells.Item(75, 8).Value = 1480     # H75 Muertes

# Row 76 - Hungria
$ws.Cells.Item(76, 2).Value = 39862    # B76 Casos totales
$ws.Cells.Item(76, 3).Value = 1025     # C76 Nuevos casos
$ws.Cells.Item(76, 4).Value = 11753    # D76 Casos activos
$ws.Cells.Item(76, 5).Value = 27113    # E76 Recuperados
$ws.Cells.Item(76, 7).Value = 28       # G76 Muertes hoy
$ws.Cells.Item(76, 8).Value = 996      # H76 Muertes

# Row 81 - El Salvador
$ws.Cells.Item(81, 2).Value = 30480    # B81 Casos totales
$ws.Cells.Item(81, 4).Value = 25857    # D81 Casos activos
$ws.Cells.Item(81, 5).Value = 3724     # E81 Recuperados
$ws.Cells.Item(81, 8).Value = 899      # H81 Muertes

# Row 101 - Georgia
$ws.Cells.Item(101, 2).Value = 12841   # B101 Casos totales
$ws.Cells.Item(101, 3).Value = 569     # C101 Nuevos casos
$ws.Cells.Item(101, 4).Value = 6867    # D101 Casos activos
$ws.Cells.Item(101, 5).Value = 5872    # E101 Recuperados
$ws.Cells.Item(101, 7).Value = 9       # G101 Muertes hoy
$ws.Cells.Item(101, 8).Value = 102     # H101 Muertes
